$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 555, shifting existing rows 555-629 down to 556-630.
# xlShiftDown = -4121
$ws.Rows.Item(555).Insert(-4121)

# Row 556 now holds what used to be row 555's data - use it as the template for
# the columns that stay constant throughout this table.
$colA = $ws.Cells.Item(556, 1).Value()
$colB = $ws.Cells.Item(556, 2).Value()
$colC = $ws.Cells.Item(556, 3).Value()
$colE = $ws.Cells.Item(556, 5).Value()
$colF = $ws.Cells.Item(556, 6).Value()
$colG = $ws.Cells.Item(556, 7).Value()
$colH = $ws.Cells.Item(556, 8).Value()
$colI = $ws.Cells.Item(556, 9).Value()
$colJ = $ws.Cells.Item(556, 10).Value()
$colK = $ws.Cells.Item(556, 11).Value()
$colL = $ws.Cells.Item(556, 12).Value()
$colQ = $ws.Cells.Item(556, 17).Value()
$colT = $ws.Cells.Item(556, 20).Value()

$ws.Cells.Item(555, 1).Value = $colA
$ws.Cells.Item(555, 2).Value = $colB
$ws.Cells.Item(555, 3).Value = $colC
$ws.Cells.Item(555, 4).Value = 45154
$ws.Cells.Item(555, 5).Value = $colE
$ws.Cells.Item(555, 6).Value = $colF
$ws.Cells.Item(555, 7).Value = $colG
$ws.Cells.Item(555, 8).Value = $colH
$ws.Cells.Item(555, 9).Value = $colI
$ws.Cells.Item(555, 10).Value = $colJ
$ws.Cells.Item(555, 11).Value = $colK
$ws.Cells.Item(555, 12).Value = $colL
$ws.Cells.Item(555, 13).Value = 185
$ws.Cells.Item(555, 14).Value = 10000
$ws.Cells.Item(555, 15).Value = 10000
$ws.Cells.Item(555, 16).Value = 10000
$ws.Cells.Item(555, 17).Value = $colQ
$ws.Cells.Item(555, 18).Value = "Brasil"
$ws.Cells.Item(555, 19).Value = 2500
$ws.Cells.Item(555, 20).Value = $colT
